$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.829.79"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.239.72"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.567"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +20.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0982"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Value = "2.571.40"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.864"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "2.236.71"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "41.823.64"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0728"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +21.61%  "
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("E38").Value = "  +14.30%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0282"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.25%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.84%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.103"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +7.54%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.08%  "
